$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 259.5122
$ws.Range("J17").Value = 263.5
$ws.Range("L17").Value = 790.5
$ws.Range("N17").Value = -1126.5
$ws.Range("H33").Value = 292.625
$ws.Range("I33").Value = 301.7931
$ws.Range("J33").Value = 204
$ws.Range("K33").Value = 301.7931
$ws.Range("L33").Value = 204
$ws.Range("M33").Value = -72.79309999999998
$ws.Range("N33").Value = -662
$ws.Range("H98").Value = 4390.5
$ws.Range("I98").Value = 4390.5
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 4390.5
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -2892.5
$ws.Range("N98").ClearContents()
$ws.Range("H118").Value = 1027
$ws.Range("I118").Value = 610
$ws.Range("K118").Value = 1830
$ws.Range("M118").Value = -173
$ws.Range("H122").Value = 4390.5
$ws.Range("I122").Value = 4390.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13171.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -10721.5
$ws.Range("N122").ClearContents()

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2102510.2
$ws.Range("I2").Value = 1916.4546
$ws.Range("K2").Value = 1916.4546
$ws.Range("M2").Value = -1803.4546
$ws.Range("H32").Value = 7014.672
$ws.Range("I32").Value = 2948.9614
$ws.Range("J32").Value = 24632.75
$ws.Range("K32").Value = 2948.9614
$ws.Range("L32").Value = 24632.75
$ws.Range("M32").Value = -2661.9614
$ws.Range("N32").Value = -25206.75
$ws.Range("H45").Value = 54352.473
$ws.Range("I45").Value = 92362.09
$ws.Range("J45").Value = 2089.25
$ws.Range("K45").Value = 92362.09
$ws.Range("L45").Value = 2089.25
$ws.Range("M45").Value = -91985.09
$ws.Range("N45").Value = -2843.25
$ws.Range("H116").Value = 2102510.2
$ws.Range("I116").Value = 1916.4546
$ws.Range("K116").Value = 1916.4546
$ws.Range("M116").Value = 377.5454
$ws.Range("H122").Value = 2444.5
$ws.Range("I122").Value = 2116
$ws.Range("J122").Value = 2700
$ws.Range("K122").Value = 6348
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -3898
$ws.Range("N122").Value = -13000
$ws.Range("H132").Value = 1920.7551
$ws.Range("I132").Value = 1676.2424
$ws.Range("J132").Value = 2425.0625
$ws.Range("K132").Value = 5028.7272
$ws.Range("L132").Value = 7275.1875
$ws.Range("M132").Value = -2498.7272
$ws.Range("N132").Value = -12335.1875

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2102510.2
$ws.Range("I3").Value = 1916.4546
$ws.Range("K3").Value = 1916.4546
$ws.Range("M3").Value = -1802.4546
$ws.Range("H86").Value = 1630.4166
$ws.Range("I86").Value = 1437.7273
$ws.Range("J86").Value = 3750
$ws.Range("K86").Value = 1437.7273
$ws.Range("L86").Value = 3750
$ws.Range("M86").Value = -314.7273
$ws.Range("N86").Value = -5996
$ws.Range("H89").Value = 1630.4166
$ws.Range("I89").Value = 1437.7273
$ws.Range("J89").Value = 3750
$ws.Range("K89").Value = 7188.636500000001
$ws.Range("L89").Value = 18750
$ws.Range("M89").Value = -1572.636500000001
$ws.Range("N89").Value = -29982

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1143.9333
$ws.Range("I16").Value = 963.875
$ws.Range("K16").Value = 963.875
$ws.Range("M16").Value = -676.875
$ws.Range("H31").Value = 1786.5483
$ws.Range("I31").Value = 1232.8667
$ws.Range("J31").Value = 2305.625
$ws.Range("K31").Value = 1232.8667
$ws.Range("L31").Value = 2305.625
$ws.Range("M31").Value = -937.8667
$ws.Range("N31").Value = -2895.625
$ws.Range("H34").Value = 1786.5483
$ws.Range("I34").Value = 1232.8667
$ws.Range("J34").Value = 2305.625
$ws.Range("K34").Value = 1232.8667
$ws.Range("L34").Value = 2305.625
$ws.Range("M34").Value = -1030.8667
$ws.Range("N34").Value = -2709.625
$ws.Range("H105").Value = 1017.6667
$ws.Range("I105").Value = 950
$ws.Range("J105").Value = 1254.5
$ws.Range("K105").Value = 950
$ws.Range("L105").Value = 1254.5
$ws.Range("M105").Value = 797
$ws.Range("N105").Value = -4748.5
$ws.Range("H113").Value = 1143.9333
$ws.Range("I113").Value = 963.875
$ws.Range("K113").Value = 963.875
$ws.Range("M113").Value = 1206.125

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 906.8642
$ws.Range("I68").Value = 615.54285
$ws.Range("J68").Value = 1128.5217
$ws.Range("K68").Value = 1846.62855
$ws.Range("L68").Value = 3385.5651
$ws.Range("M68").Value = -1035.62855
$ws.Range("N68").Value = -5007.5651
$ws.Range("H71").Value = 906.8642
$ws.Range("I71").Value = 615.54285
$ws.Range("J71").Value = 1128.5217
$ws.Range("K71").Value = 5539.88565
$ws.Range("L71").Value = 10156.6953
$ws.Range("M71").Value = -1483.88565
$ws.Range("N71").Value = -18268.6953
$ws.Range("H131").Value = 1283758
$ws.Range("I131").Value = 2782.5
$ws.Range("J131").Value = 1352999.9
$ws.Range("K131").Value = 8347.5
$ws.Range("L131").Value = 4058999.7
$ws.Range("M131").Value = -3307.5
$ws.Range("N131").Value = -4069079.7

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1285
$ws.Range("I97").Value = 1285
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1285
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -789
$ws.Range("N97").ClearContents()
$ws.Range("H122").Value = 10811.889
$ws.Range("I122").Value = 14051.167
$ws.Range("J122").Value = 4333.3335
$ws.Range("K122").Value = 42153.501
$ws.Range("L122").Value = 13000.0005
$ws.Range("M122").Value = -39703.501
$ws.Range("N122").Value = -17900.0005

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1975.375
$ws.Range("I40").Value = 1966.3334
$ws.Range("J40").Value = 2002.5
$ws.Range("K40").Value = 1966.3334
$ws.Range("L40").Value = 2002.5
$ws.Range("M40").Value = -1830.3334
$ws.Range("N40").Value = -2274.5
$ws.Range("H132").Value = 5478.5674
$ws.Range("I132").Value = 5949.1665
$ws.Range("K132").Value = 17847.4995
$ws.Range("M132").Value = -15317.4995

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 36000
$ws.Range("J75").Value = 36000
$ws.Range("L75").Value = 36000
$ws.Range("N75").Value = -37872
$ws.Range("H78").Value = 36000
$ws.Range("J78").Value = 36000
$ws.Range("L78").Value = 108000
$ws.Range("N78").Value = -117360
$ws.Range("H136").Value = 30558.912
$ws.Range("I136").Value = 835.64703
$ws.Range("K136").Value = 2506.94109
$ws.Range("M136").Value = 43.0589100000002
